$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '87.169.95'
Set-TextValue 'E2' '  -2.47%  '

# Row 3
Set-TextValue 'D3' '3.139.08'
Set-TextValue 'E3' '  -6.65%  '

# Row 4
Set-TextValue 'E4' '  -0.03%  '

# Row 5
Set-TextValue 'D5' '202.39'
Set-TextValue 'E5' '  -7.89%  '

# Row 6
Set-TextValue 'D6' '604.28'
Set-TextValue 'E6' '  -6.87%  '

# Row 7
Set-TextValue 'D7' '0.372'
Set-TextValue 'E7' '  -9.02%  '

# Row 8
Set-TextValue 'D8' '0.655'
Set-TextValue 'E8' '  +6.44%  '

# Row 9
Set-TextValue 'E9' '  -0.01%  '

# Row 10
Set-TextValue 'D10' '3.131.52'

# Row 11
Set-TextValue 'D11' '0.527'
Set-TextValue 'E11' '  -13.94%  '

# Row 12
Set-TextValue 'E12' '  +4.86%  '

# Row 13
Set-TextValue 'D13' '0.0000240'
Set-TextValue 'E13' '  -16.82%  '

# Row 14
Set-TextValue 'D14' '3.720.62'
Set-TextValue 'E14' '  -6.64%  '

# Row 15
Set-TextValue 'D15' '5.20'
Set-TextValue 'E15' '  -6.53%  '

# Row 16
Set-TextValue 'D16' '86.891.38'
Set-TextValue 'E16' '  -2.67%  '

# Row 17
Set-TextValue 'D17' '31.69'
Set-TextValue 'E17' '  -13.30%  '

# Row 18
Set-TextValue 'D18' '3.155.53'
Set-TextValue 'E18' '  -6.05%  '

# Row 19
Set-TextValue 'D19' '2.98'
Set-TextValue 'E19' '  -5.22%  '

# Row 20
Set-TextValue 'D20' '13.25'
Set-TextValue 'E20' '  -10.83%  '

# Row 21
Set-TextValue 'D21' '410.97'
Set-TextValue 'E21' '  -10.28%  '

# Row 22
Set-TextValue 'D22' '8.40'
Set-TextValue 'E22' '  -13.36%  '

# Row 23
Set-TextValue 'D23' '5.04'
Set-TextValue 'E23' '  -9.36%  '

# Row 24
Set-TextValue 'D24' '5.08'
Set-TextValue 'E24' '  -8.68%  '

# Row 25
Set-TextValue 'D25' '11.82'
Set-TextValue 'E25' '  -7.95%  '

# Row 26
Set-TextValue 'D26' '3.308.74'
Set-TextValue 'E26' '  -6.12%  '

# Row 27
Set-TextValue 'D27' '73.03'
Set-TextValue 'E27' '  -7.82%  '

# Row 28
Set-TextValue 'D28' '0.0000128'
Set-TextValue 'E28' '  -10.54%  '

# Row 29
Set-TextValue 'D29' '0.999'
Set-TextValue 'E29' '  -0.04%  '

# Row 30
Set-TextValue 'D30' '0.160'
Set-TextValue 'E30' '  -19.82%  '

# Row 31
Set-TextValue 'E31' '  +0.14%  '

# Row 32
Set-TextValue 'D32' '532.93'
Set-TextValue 'E32' '  -10.65%  '

# Row 33
Set-TextValue 'D33' '8.19'
Set-TextValue 'E33' '  -13.34%  '

# Row 34
Set-TextValue 'E34' '  -17.97%  '

# Row 35
Set-TextValue 'D35' '1.82'
Set-TextValue 'E35' '  -13.69%  '

# Row 36
Set-TextValue 'D36' '6.58'
Set-TextValue 'E36' '  -10.77%  '

# Row 37
Set-TextValue 'E37' '  -9.00%  '

# Row 38
Set-TextValue 'D38' '21.71'
Set-TextValue 'E38' '  -7.74%  '

# Row 39
Set-TextValue 'D39' '21.76'
Set-TextValue 'E39' '  -0.39%  '

# Row 40
Set-TextValue 'D40' '0.997'
Set-TextValue 'E40' '  -0.11%  '

# Row 41
Set-TextValue 'D41' '2.95'
Set-TextValue 'E41' '  -8.34%  '

# Row 42
Set-TextValue 'E42' '  +0.03%  '

# Row 43
Set-TextValue 'B43' 'Stacks'
Set-TextValue 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D43' '1.87'
Set-TextValue 'E43' '  -13.38%  '

# Row 44
Set-TextValue 'B44' 'PolygonEcosystemToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D44' '0.367'
Set-TextValue 'E44' '  -14.64%  '

# Row 45
Set-TextValue 'D45' '147.41'
Set-TextValue 'E45' '  -6.72%  '

# Row 46
Set-TextValue 'D46' '170.92'
Set-TextValue 'E46' '  -10.26%  '

# Row 47
Set-TextValue 'E47' '  -7.35%  '

# Row 48
Set-TextValue 'E48' '  +4.83%  '

# Row 49
Set-TextValue 'D49' '1.24'
Set-TextValue 'E49' '  -16.27%  '

# Row 50
Set-TextValue 'D50' '3.93'
Set-TextValue 'E50' '  -12.90%  '

# Row 51
Set-TextValue 'D51' '0.688'
Set-TextValue 'E51' '  -12.53%  '
